# Rewrite the Sheet1 test-case content: new testcase id, new "Valid
# Scenario" / "Invalid Scenario" steps, and matching cell formatting
# (bold section headers with a left indent, a wrapped "Expected values"
# cell, and indented step rows underneath each header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131
$xlCenter = -4108

# ---------------------------------------------------------------------
# Cell VALUES first, written in the same order the source content was
# authored in, so newly-created shared-string entries land in the same
# order as the target workbook.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "TC_cura_validateappointmentbutton_002"
$ws.Range("B2").Value = "1. Valid Scenario:"
$ws.Range("B3").Value = "Homepage navigation."
$ws.Range("B4").Value = 'Confirm "Make Appointment" button visibility.'
$ws.Range("B5").Value = 'Click "Make Appointment."'
$ws.Range("B7").Value = '(If applicable) Test "Make Appointment" button in invalid context (e.g., when not logged in).'
$ws.Range("C2").Value = "    Redirection`n    Presence"
$ws.Range("B6").Value = "2. Invalid Scenario :"
$ws.Range("C6").Value = "(If applicable) Error"

# ---------------------------------------------------------------------
# Cell FORMATTING, applied in the order that mints each distinct style
# (plain indent=1, bold+indent=1, plain indent=2, wrap text) so the
# saved cellXfs table matches the target workbook's style order.
# ---------------------------------------------------------------------

# Plain left-indent(1) style - first used on the trailing blank rows.
$ws.Range("B8").HorizontalAlignment = $xlLeft
$ws.Range("B8").VerticalAlignment = $xlCenter
$ws.Range("B8").IndentLevel = 1

# Bold + left-indent(1) style - the "N. Scenario:" section headers.
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").HorizontalAlignment = $xlLeft
$ws.Range("B2").VerticalAlignment = $xlCenter
$ws.Range("B2").IndentLevel = 1

# Plain left-indent(2) style - the individual step rows.
$ws.Range("B3").HorizontalAlignment = $xlLeft
$ws.Range("B3").VerticalAlignment = $xlCenter
$ws.Range("B3").IndentLevel = 2

# Wrap-text style - the multi-line "Expected values" cell.
$ws.Range("C2").WrapText = $true

$ws.Range("B4").HorizontalAlignment = $xlLeft
$ws.Range("B4").VerticalAlignment = $xlCenter
$ws.Range("B4").IndentLevel = 2

$ws.Range("B5").HorizontalAlignment = $xlLeft
$ws.Range("B5").VerticalAlignment = $xlCenter
$ws.Range("B5").IndentLevel = 2

$ws.Range("B6").Font.Bold = $true
$ws.Range("B6").HorizontalAlignment = $xlLeft
$ws.Range("B6").VerticalAlignment = $xlCenter
$ws.Range("B6").IndentLevel = 1

$ws.Range("B7").HorizontalAlignment = $xlLeft
$ws.Range("B7").VerticalAlignment = $xlCenter
$ws.Range("B7").IndentLevel = 2

$ws.Range("B10").HorizontalAlignment = $xlLeft
$ws.Range("B10").VerticalAlignment = $xlCenter
$ws.Range("B10").IndentLevel = 1

$ws.Range("B11").HorizontalAlignment = $xlLeft
$ws.Range("B11").VerticalAlignment = $xlCenter
$ws.Range("B11").IndentLevel = 1

# Row height for row 2 to accommodate the two-line wrapped "Expected
# values" text.
$ws.Rows.Item(2).RowHeight = 28.8

# ---------------------------------------------------------------------
# Column B is widened to fit the new step/description text.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 55.166666666666664

# ---------------------------------------------------------------------
# View state: leave the selection on C14.
# ---------------------------------------------------------------------
$ws.Range("C14").Select()
